$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 582, shifting existing rows 582:672 down to 585:675.
$ws.Rows("582:584").Insert()

# Common (constant across the whole dataset) values for the new rows.
$commonA = 7
$commonB = "Terminal Hortofrutícola Agro Chillán"
$commonC = "Ñuble"
$commonE = 16
$commonF = 100114001
$commonG = "Papa"
$commonN = "`$/saco 25 kilos"
$commonQ = 25
$commonR = "Hortaliza"

# New row 582
$ws.Range("A582").Value = $commonA
$ws.Range("B582").Value = $commonB
$ws.Range("C582").Value = $commonC
$ws.Range("D582").Value = 45127
$ws.Range("E582").Value = $commonE
$ws.Range("F582").Value = $commonF
$ws.Range("G582").Value = $commonG
$ws.Range("H582").Value = "Asterix"
$ws.Range("I582").Value = "1a (guarda)"
$ws.Range("J582").Value = 150
$ws.Range("K582").Value = 18000
$ws.Range("L582").Value = 18000
$ws.Range("M582").Value = 18000
$ws.Range("N582").Value = $commonN
$ws.Range("O582").Value = "Región de Los Lagos"
$ws.Range("P582").Value = 720
$ws.Range("Q582").Value = $commonQ
$ws.Range("R582").Value = $commonR

# New row 583
$ws.Range("A583").Value = $commonA
$ws.Range("B583").Value = $commonB
$ws.Range("C583").Value = $commonC
$ws.Range("D583").Value = 45127
$ws.Range("E583").Value = $commonE
$ws.Range("F583").Value = $commonF
$ws.Range("G583").Value = $commonG
$ws.Range("H583").Value = "Asterix"
$ws.Range("I583").Value = "1a (guarda)"
$ws.Range("J583").Value = 100
$ws.Range("K583").Value = 17000
$ws.Range("L583").Value = 17000
$ws.Range("M583").Value = 17000
$ws.Range("N583").Value = $commonN
$ws.Range("O583").Value = "Región de Ñuble"
$ws.Range("P583").Value = 680
$ws.Range("Q583").Value = $commonQ
$ws.Range("R583").Value = $commonR

# New row 584
$ws.Range("A584").Value = $commonA
$ws.Range("B584").Value = $commonB
$ws.Range("C584").Value = $commonC
$ws.Range("D584").Value = 45127
$ws.Range("E584").Value = $commonE
$ws.Range("F584").Value = $commonF
$ws.Range("G584").Value = $commonG
$ws.Range("H584").Value = "Asterix"
$ws.Range("I584").Value = "2a (guarda)"
$ws.Range("J584").Value = 180
$ws.Range("K584").Value = 16000
$ws.Range("L584").Value = 16000
$ws.Range("M584").Value = 16000
$ws.Range("N584").Value = $commonN
$ws.Range("O584").Value = "Región de Los Lagos"
$ws.Range("P584").Value = 640
$ws.Range("Q584").Value = $commonQ
$ws.Range("R584").Value = $commonR
